$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("index")

# Update the translation value for the "nav_servicio" key from "Servicio" to "Servicios"
$ws.Range("B4").Value = "Servicios"

$wb.Save()
